# Score_Manager.xlsx update:
#  - Database Updated: add a "CPU" sub-score row under the "Processor" group
#    (G5/H5/I5 were blank placeholder cells; Excel's "Prioritise" script now
#    fills them in with Name=CPU, Score=5, Rating=M).
#  - Debug Data updated: the active selection moved to H6 (the cell right
#    below the newly entered data) as a result of the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "CPU"
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = "M"

$ws.Range("H6").Select() | Out-Null
